$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.267.03'
$ws.Range('E2').Value = '  +1.38%  '
$ws.Range('D3').Value = '2.423.15'
$ws.Range('E3').Value = '  +1.93%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '562.80'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.16'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.14%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +1.65%  '
$ws.Range('D9').Value = '2.420.66'
$ws.Range('E9').Value = '  +1.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.110'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.75%  '
$ws.Range('E11').Value = '  -2.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.38'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.53%  '
$ws.Range('E13').Value = '  +0.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.92'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.62%  '
$ws.Range('E15').Value = '  +3.45%  '
$ws.Range('D16').Value = '2.860.53'
$ws.Range('E16').Value = '  +1.95%  '
$ws.Range('D17').Value = '62.107.07'
$ws.Range('E17').Value = '  +1.31%  '
$ws.Range('D18').Value = '2.419.99'
$ws.Range('E18').Value = '  +1.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.32'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.07%  '
$ws.Range('E20').Value = '  +1.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '324.29'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.76'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.96%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.59'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.71'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.74%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.93'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.93%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '585.51'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Value = '2.541.47'
$ws.Range('E28').Value = '  +1.96%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').Value = '0.0₃0946'
$ws.Range('E30').Value = '  +4.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.47'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.19%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.25'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.13%  '
$ws.Range('E33').Value = '  +1.45%  '
$ws.Range('E34').Value = '  +2.52%  '
$ws.Range('E35').Value = '  +1.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.73'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.56%  '
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.80'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.385'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '152.71'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.69'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.07%  '
$ws.Range('E42').Value = '  -3.60%  '
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.33'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +8.62%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '150.34'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.47%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.67'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.83%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0538'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.39'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.595'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0925'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.89%  '
$ws.Range('E51').Value = '  +2.22%  '
